# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same column layout: 基金代码/基金名称/
#    基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名) right after it,
#    rename the copy to "2022-Q1", and overwrite its data with the new
#    quarter's fund-holding figures (11 funds).
# 2. Prepend a "2022-Q1" row to the "总计" (totals) summary sheet and
#    renumber the index column.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    # Force the cell to be stored as text even when the value looks
    # numeric (e.g. "33.17" or a fund code with a leading zero like
    # "010190"), matching the source data's inline-string typing.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Set-NumValue($cell, $val) {
    # Plain numeric cell, reset to the default (unstyled) look - undoes
    # any formatting a row Insert() may have carried over from a
    # neighbouring row.
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" holdings sheet from the "2021-Q4" template
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)
$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

# Trim the copied sheet down from 20 data rows to the 11 rows of the new
# quarter (rows 2..12 stay, rows 13..20 go).
$ws.Rows("13:20").Delete()

$fundRows = @(
    @("010190", "嘉实价值发现三个月定期开放混合", "33.17", "93.94", "4.41", "1.4628", 9),
    @("000480", "东方红新动力灵活配置混合",       "15.38", "72.90", "7.11", "1.0935", 1),
    @("001044", "嘉实新消费股票",                 "8.92",  "80.25", "5.39", "0.4808", 7),
    @("001564", "东方红京东大数据灵活配置混合",    "11.40", "69.58", "3.91", "0.4457", 4),
    @("004355", "嘉实丰和灵活配置混合",            "9.22",  "85.59", "4.50", "0.4149", 7),
    @("008704", "广发高股息优享混合A",             "3.52",  "94.09", "9.20", "0.3238", 2),
    @("004119", "广发创新驱动灵活配置混合",        "2.79",  "94.36", "9.71", "0.2709", 1),
    @("008705", "广发高股息优享混合C",             "0.82",  "94.09", "9.20", "0.0754", 2),
    @("519959", "长信多利灵活配置混合",            "1.45",  "85.11", "3.73", "0.0541", 10),
    @("013488", "长信多利灵活配置混合D",           "1.45",  "85.11", "3.73", "0.0541", 10),
    @("519987", "长信恒利优势混合",                "0.22",  "82.39", "3.94", "0.0087", 10)
)

$r = 2
foreach ($row in $fundRows) {
    $ws.Cells.Item($r, 1).Value = ($r - 2)
    Set-TextValue $ws.Cells.Item($r, 2) $row[0]
    Set-TextValue $ws.Cells.Item($r, 3) $row[1]
    Set-TextValue $ws.Cells.Item($r, 4) $row[2]
    Set-TextValue $ws.Cells.Item($r, 5) $row[3]
    Set-TextValue $ws.Cells.Item($r, 6) $row[4]
    Set-TextValue $ws.Cells.Item($r, 7) $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 2: prepend the new quarter to the "总计" totals sheet
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows("2:2").Insert()

Set-NumValue $tot.Cells.Item(2, 1) 0
Set-TextValue $tot.Cells.Item(2, 2) "2022-Q1"
Set-NumValue $tot.Cells.Item(2, 3) 11
Set-NumValue $tot.Cells.Item(2, 4) 4.68

# The freshly-inserted row picked up formatting from its neighbour; put
# the A2 "index column" look (bold + thin border, same as A3:A7) back by
# copying the format from an already-correct cell in that column.
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

$totRows = @(
    @("2021-Q4", 19, 5.52),
    @("2021-Q3", 20, 5.41),
    @("2021-Q2", 35, 7.85),
    @("2021-Q1", 12, 2.61),
    @("2020-Q4", 11, 6.77)
)

$r = 3
foreach ($row in $totRows) {
    $tot.Cells.Item($r, 1).Value = ($r - 2)
    Set-TextValue $tot.Cells.Item($r, 2) $row[0]
    $tot.Cells.Item($r, 3).Value = $row[1]
    $tot.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}
